# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = "Chainlink"
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("B44").Value = "Maker"
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("B48").Value = "FTXToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D2").Value = "37.675.11"
$ws.Range("D3").Value = "2.075.25"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "58.14"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.389"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0777"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.84"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "2.381.19"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.765"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "2.069.31"
$ws.Range("D18").Value = "37.598.98"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "0.0₃0828"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.58"
$ws.Range("D32").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.52"
$ws.Range("D35").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.39"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0958"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.91"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "1.486.39"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.72"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").Value = "2.265.67"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -1.47%  "
$ws.Range("E9").Value = "  -1.69%  "
$ws.Range("E10").Value = "  -1.59%  "
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("E15").Value = "  -1.22%  "
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("E21").Value = "  -2.52%  "
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("E25").Value = "  -2.74%  "
$ws.Range("E26").Value = "  +2.69%  "
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("E28").Value = "  -5.09%  "
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("E30").Value = "  -4.23%  "
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("E32").Value = "  -3.34%  "
$ws.Range("E33").Value = "  -1.06%  "
$ws.Range("E34").Value = "  -0.72%  "
$ws.Range("E35").Value = "  +0.99%  "
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("E37").Value = "  -3.68%  "
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("E39").Value = "  -2.08%  "
$ws.Range("E40").Value = "  +4.00%  "
$ws.Range("E41").Value = "  -0.86%  "
$ws.Range("E42").Value = "  -2.46%  "
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("E44").Value = "  +2.52%  "
$ws.Range("E45").Value = "  +2.96%  "
$ws.Range("E46").Value = "  -2.80%  "
$ws.Range("E47").Value = "  -2.26%  "
$ws.Range("E48").Value = "  -3.46%  "
$ws.Range("E49").Value = "  -1.44%  "
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("E51").Value = "  -0.52%  "
